# --- Set up workbook/sheet structure ---
$wb = $excel.ActiveWorkbook
$general = $wb.Worksheets.Item("general")

# Move selection on the general sheet to A2 (matches target diff)
$general.Range("A2").Select()

# Insert the new sheet right after "general"
$ws = $wb.Worksheets.Add($null, $general)
$ws.Name = "Sheet1"

# --- Header row (row 1) ---
$ws.Range("A1").Value = 'TPN'
$ws.Range("B1").Value = 'Description'
$ws.Range("C1").Value = 'Color'
$ws.Range("D1").Value = 'Configuration'
$ws.Range("E1").Value = 'Lens Color'
$ws.Range("F1").Value = 'Lens Transparency'
$ws.Range("G1").Value = 'Millicandela Rating'
$ws.Range("H1").Value = 'Lens Style'
$ws.Range("I1").Value = 'Lens Size'
$ws.Range("J1").Value = 'Voltage - Forward (Vf) (Typ)'
$ws.Range("K1").Value = 'Current - Test'
$ws.Range("L1").Value = 'Viewing Angle'
$ws.Range("M1").Value = 'Mounting Type'
$ws.Range("N1").Value = 'Wavelength - Dominant'
$ws.Range("O1").Value = 'Wavelength - Peak'
$ws.Range("P1").Value = 'Features'
$ws.Range("Q1").Value = 'Package / Case'
$ws.Range("R1").Value = 'Supplier Device Package'
$ws.Range("S1").Value = 'Size / Dimension'
$ws.Range("T1").Value = 'Height (Max)'
$ws.Range("U1").Value = 'Library Ref'
$ws.Range("V1").Value = 'Footprint Ref'
$ws.Range("W1").Value = 'Manufacturer 1'
$ws.Range("X1").Value = 'Manufacturer 1 PN'
$ws.Range("Y1").Value = 'Supplier 1'
$ws.Range("Z1").Value = 'Supplier 1 PN'
$ws.Range("AA1").Value = 'Supplier 1 Link'
$ws.Range("AB1").Value = 'Manufacturer 2'
$ws.Range("AC1").Value = 'Manufacturer 2 PN'
$ws.Range("AD1").Value = 'Supplier 2'
$ws.Range("AE1").Value = 'Supplier 2 PN'

# Header row formatting: bold (matches style index 1) + text format on C1 (style index 2)
$ws.Range("A1:AE1").Font.Bold = $true
$ws.Range("C1").NumberFormat = "@"

# --- Row 2: LED RED DIFFUSED T-1 3/4 T/H ---
$ws.Range("A2").Formula = '="LED-"&TEXT(ROW()-1,"000000")'
$ws.Range("B2").Value = 'LED RED DIFFUSED T-1 3/4 T/H'
$ws.Range("C2").Value = 'Red'
$ws.Range("D2").Value = 'Standard'
$ws.Range("E2").Value = 'Red'
$ws.Range("F2").Value = 'Diffused'
$ws.Range("G2").Value = '19mcd'
$ws.Range("H2").Value = 'Round with Domed Top'
$ws.Range("I2").Value = '5mm, T-1 3/4'
$ws.Range("J2").Value = '2V'
$ws.Range("K2").Value = '10mA'
$ws.Range("L2").Value = '36°'
$ws.Range("M2").Value = 'Through Hole'
$ws.Range("N2").Value = '623nm'
$ws.Range("O2").Value = '635nm'
$ws.Range("P2").Value = '-'
$ws.Range("Q2").Value = 'Radial'
$ws.Range("R2").Value = 'T-1 3/4'
$ws.Range("S2").Value = '-'
$ws.Range("T2").Value = '8.60mm'
$ws.Range("U2").Value = 'LED-Red'
$ws.Range("W2").Value = 'Lite-On Inc.'
$ws.Range("X2").Value = 'LTL-4223'
$ws.Range("Y2").Value = 'Digi-Key'
$ws.Range("Z2").Value = '160-1127-ND'
$ws.Range("AA2").Value = 'https://www.digikey.com/product-detail/en/lite-on-inc/LTL-4223/160-1127-ND/200395'

# --- Row 3: LED BLUE CLEAR 5MM ROUND T/H ---
$ws.Range("A3:A4").Formula = '="LED-"&TEXT(ROW()-1,"000000")'
$ws.Range("B3").Value = 'LED BLUE CLEAR 5MM ROUND T/H'
$ws.Range("C3").Value = 'Blue'
$ws.Range("D3").Value = 'Standard'
$ws.Range("E3").Value = 'Colorless'
$ws.Range("F3").Value = 'Clear'
$ws.Range("G3").Value = '7065mcd'
$ws.Range("H3").Value = 'Round with Domed Top'
$ws.Range("I3").Value = '5.0mm Dia'
$ws.Range("J3").Value = '3.2V'
$ws.Range("K3").Value = '20mA'
$ws.Range("L3").Value = '30°'
$ws.Range("M3").Value = 'Through Hole'
$ws.Range("N3").Value = '470nm'
$ws.Range("O3").Value = '-'
$ws.Range("P3").Value = '-'
$ws.Range("Q3").Value = 'Radial'
$ws.Range("R3").Value = '5-mm Round'
$ws.Range("S3").Value = '-'
$ws.Range("T3").Value = '8.60mm'
$ws.Range("U3").Value = 'LED-Blue'
$ws.Range("W3").Value = 'Cree Inc.'
$ws.Range("X3").Value = 'C503B-BCN-CV0Z0461'
$ws.Range("Y3").Value = 'Digi-Key'
$ws.Range("Z3").Value = 'C503B-BCN-CV0Z0461-ND'
$ws.Range("AA3").Value = 'https://www.digikey.com/product-detail/en/cree-inc/C503B-BCN-CV0Z0461/C503B-BCN-CV0Z0461-ND/1922945'

# --- Row 4: LED BLUE CLEAR T-1 3/4 T/H ---
$ws.Range("B4").Value = 'LED BLUE CLEAR T-1 3/4 T/H'
$ws.Range("C4").Value = 'Blue'
$ws.Range("D4").Value = 'Standard'
$ws.Range("E4").Value = 'Colorless'
$ws.Range("F4").Value = 'Clear'
$ws.Range("G4").Value = '1200mcd'
$ws.Range("H4").Value = 'Round with Domed Top'
$ws.Range("I4").Value = '5mm, T-1 3/4'
$ws.Range("J4").Value = '3.5V'
$ws.Range("K4").Value = '10mA'
$ws.Range("L4").Value = '36°'
$ws.Range("M4").Value = 'Through Hole'
$ws.Range("N4").Value = '470nm'
$ws.Range("O4").Value = '468nm'
$ws.Range("P4").Value = '-'
$ws.Range("Q4").Value = 'Radial'
$ws.Range("R4").Value = 'T-1 3/4'
$ws.Range("S4").Value = '-'
$ws.Range("T4").Value = '8.70mm'
$ws.Range("U4").Value = 'LED-Blue'
$ws.Range("W4").Value = 'Lite-On Inc.'
$ws.Range("X4").Value = 'LTL2T3TBK5'
$ws.Range("Y4").Value = 'Digi-Key'
$ws.Range("Z4").Value = '160-1610-ND'
$ws.Range("AA4").Value = 'https://www.digikey.com/product-detail/en/lite-on-inc/LTL2T3TBK5/160-1610-ND/573515'

# --- Page setup (matches target diff: portrait orientation) ---
$ws.PageSetup.Orientation = 1

# --- Final view state on the new sheet (matches target diff) ---
$ws.Range("AC4").Select()
